$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 204, pushing existing rows 204..270 down to 205..271
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44663
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100112040
$ws.Range("G204").Value = "Cilantro"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 250
$ws.Range("K204").Value = 10000
$ws.Range("L204").Value = 10000
$ws.Range("M204").Value = 10000
$ws.Range("N204").Value = "`$/caja 36 atados"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 278
$ws.Range("Q204").Value = 36
$ws.Range("R204").Value = "Hortaliza"
